$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order id timestamps refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555138297467"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555155708125"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555155718138"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555156348152"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555157128134"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555137917488.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555138117464.csv"
$ws1.Range("B4").Value = "go_stims-1651255513813748.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255513827749.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_6-16512555140874515.csv"
$ws2.Range("B3").Value = "TB-16512555153964458.csv"
$ws2.Range("B4").Value = "OB-16512555145544455.csv"
$ws2.Range("B5").Value = "TB-16512555150674524.csv"
$ws2.Range("B6").Value = "ZB-match_8-16512555138357472.csv"
$ws2.Range("B7").Value = "ZB-match_9-16512555141354463.csv"
$ws2.Range("B8").Value = "TB-1651255515553813.csv"
$ws2.Range("B9").Value = "OB-16512555143534462.csv"
$ws2.Range("B10").Value = "OB-1651255514288447.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555156028135.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555155788143.csv"
$ws4.Range("B4").Value = "MM_stims-16512555156178145.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555156038141.csv"
$ws4.Range("B6").Value = "MM_stims-16512555156338131.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555156188154.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512555156808136.csv"
$ws5.Range("B3").Value = "SAT_stims-1651255515665813.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555156408129.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555156968114.csv"
